$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

$ws.Cells.Item(11, 4).Value = 3622668
$ws.Cells.Item(11, 5).Value = 4070396
$ws.Cells.Item(11, 6).Value = 4350370
$ws.Cells.Item(11, 7).Value = 4307063
$ws.Cells.Item(11, 8).Value = 5043948
$ws.Cells.Item(11, 9).Value = 4950335
$ws.Cells.Item(11, 10).Value = 6973055
$ws.Cells.Item(11, 11).Value = 8037511
$ws.Cells.Item(11, 12).Value = 9436084
$ws.Cells.Item(11, 13).Value = 7896746

$ws.Cells.Item(12, 4).Value = -3249932
$ws.Cells.Item(12, 5).Value = -3750008
$ws.Cells.Item(12, 6).Value = -3938742
$ws.Cells.Item(12, 7).Value = -3861092
$ws.Cells.Item(12, 8).Value = -4503178
$ws.Cells.Item(12, 9).Value = -4550403
$ws.Cells.Item(12, 10).Value = -6308166
$ws.Cells.Item(12, 11).Value = -6737622
$ws.Cells.Item(12, 12).Value = -8629979
$ws.Cells.Item(12, 13).Value = -7056472

$ws.Cells.Item(13, 4).Value = 372736
$ws.Cells.Item(13, 5).Value = 320388
$ws.Cells.Item(13, 6).Value = 411628
$ws.Cells.Item(13, 7).Value = 445971
$ws.Cells.Item(13, 8).Value = 540770
$ws.Cells.Item(13, 9).Value = 399932
$ws.Cells.Item(13, 10).Value = 664889
$ws.Cells.Item(13, 11).Value = 1299889
$ws.Cells.Item(13, 12).Value = 806105
$ws.Cells.Item(13, 13).Value = 840274

$ws.Cells.Item(14, 4).Value = -156664
$ws.Cells.Item(14, 5).Value = -94678
$ws.Cells.Item(14, 6).Value = -191341
$ws.Cells.Item(14, 7).Value = -251825
$ws.Cells.Item(14, 8).Value = -226778
$ws.Cells.Item(14, 9).Value = -109457
$ws.Cells.Item(14, 10).Value = -285578
$ws.Cells.Item(14, 11).Value = -258374
$ws.Cells.Item(14, 12).Value = -308155
$ws.Cells.Item(14, 13).Value = -383299

$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 13).Value = 0

$ws.Cells.Item(16, 4).Value = 1469
$ws.Cells.Item(16, 5).Value = -6614
$ws.Cells.Item(16, 6).Value = 9978
$ws.Cells.Item(16, 7).Value = 3857
$ws.Cells.Item(16, 8).Value = 2080
$ws.Cells.Item(16, 9).Value = 2237
$ws.Cells.Item(16, 10).Value = 908
$ws.Cells.Item(16, 11).Value = 0
$ws.Cells.Item(16, 12).Value = 0
$ws.Cells.Item(16, 13).Value = 445

$ws.Cells.Item(17, 4).Value = 217541
$ws.Cells.Item(17, 5).Value = 219096
$ws.Cells.Item(17, 6).Value = 230265
$ws.Cells.Item(17, 7).Value = 198003
$ws.Cells.Item(17, 8).Value = 316072
$ws.Cells.Item(17, 9).Value = 292712
$ws.Cells.Item(17, 10).Value = 380219
$ws.Cells.Item(17, 11).Value = 1041515
$ws.Cells.Item(17, 12).Value = 497950
$ws.Cells.Item(17, 13).Value = 457420

$ws.Cells.Item(18, 4).Value = -19934
$ws.Cells.Item(18, 5).Value = -26355
$ws.Cells.Item(18, 6).Value = -16457
$ws.Cells.Item(18, 7).Value = -11185
$ws.Cells.Item(18, 8).Value = -4079
$ws.Cells.Item(18, 9).Value = -15861
$ws.Cells.Item(18, 10).Value = -16333
$ws.Cells.Item(18, 11).Value = -20047
$ws.Cells.Item(18, 12).Value = -6310
$ws.Cells.Item(18, 13).Value = -1234

$ws.Cells.Item(19, 4).Value = 10307
$ws.Cells.Item(19, 5).Value = 11702
$ws.Cells.Item(19, 6).Value = 10842
$ws.Cells.Item(19, 7).Value = 16977
$ws.Cells.Item(19, 8).Value = 18286
$ws.Cells.Item(19, 9).Value = 18263
$ws.Cells.Item(19, 10).Value = 16564
$ws.Cells.Item(19, 11).Value = 38300
$ws.Cells.Item(19, 12).Value = 31325
$ws.Cells.Item(19, 13).Value = 67603

$ws.Cells.Item(20, 4).Value = 207914
$ws.Cells.Item(20, 5).Value = 204443
$ws.Cells.Item(20, 6).Value = 224650
$ws.Cells.Item(20, 7).Value = 203795
$ws.Cells.Item(20, 8).Value = 330279
$ws.Cells.Item(20, 9).Value = 295114
$ws.Cells.Item(20, 10).Value = 380450
$ws.Cells.Item(20, 11).Value = 1059768
$ws.Cells.Item(20, 12).Value = 522965
$ws.Cells.Item(20, 13).Value = 523789

$ws.Cells.Item(21, 4).Value = -46712
$ws.Cells.Item(21, 5).Value = -53275
$ws.Cells.Item(21, 6).Value = -50403
$ws.Cells.Item(21, 7).Value = -45481
$ws.Cells.Item(21, 8).Value = -74031
$ws.Cells.Item(21, 9).Value = -22853
$ws.Cells.Item(21, 10).Value = -84978
$ws.Cells.Item(21, 11).Value = -237114
$ws.Cells.Item(21, 12).Value = -116945
$ws.Cells.Item(21, 13).Value = 38727

$ws.Cells.Item(22, 4).Value = 161202
$ws.Cells.Item(22, 5).Value = 151168
$ws.Cells.Item(22, 6).Value = 174247
$ws.Cells.Item(22, 7).Value = 158314
$ws.Cells.Item(22, 8).Value = 256248
$ws.Cells.Item(22, 9).Value = 272261
$ws.Cells.Item(22, 10).Value = 295472
$ws.Cells.Item(22, 11).Value = 822654
$ws.Cells.Item(22, 12).Value = 406020
$ws.Cells.Item(22, 13).Value = 562516

$ws.Cells.Item(23, 4).Value = 0
$ws.Cells.Item(23, 5).Value = 0
$ws.Cells.Item(23, 6).Value = 0
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 12).Value = 0
$ws.Cells.Item(23, 13).Value = 0

$ws.Cells.Item(24, 4).Value = 161202
$ws.Cells.Item(24, 5).Value = 151168
$ws.Cells.Item(24, 6).Value = 174247
$ws.Cells.Item(24, 7).Value = 158314
$ws.Cells.Item(24, 8).Value = 256248
$ws.Cells.Item(24, 9).Value = 272261
$ws.Cells.Item(24, 10).Value = 295472
$ws.Cells.Item(24, 11).Value = 822654
$ws.Cells.Item(24, 12).Value = 406020
$ws.Cells.Item(24, 13).Value = 562516

$ws.Cells.Item(25, 4).Value = 193
$ws.Cells.Item(25, 5).Value = 181
$ws.Cells.Item(25, 6).Value = 208
$ws.Cells.Item(25, 7).Value = 189
$ws.Cells.Item(25, 8).Value = 307
$ws.Cells.Item(25, 9).Value = 326
$ws.Cells.Item(25, 10).Value = 354
$ws.Cells.Item(25, 11).Value = 984
$ws.Cells.Item(25, 12).Value = 486
$ws.Cells.Item(25, 13).Value = 487

$ws.Cells.Item(26, 4).Value = 835820
$ws.Cells.Item(26, 5).Value = 835820
$ws.Cells.Item(26, 6).Value = 835821
$ws.Cells.Item(26, 7).Value = 835820
$ws.Cells.Item(26, 8).Value = 835821
$ws.Cells.Item(26, 9).Value = 835821
$ws.Cells.Item(26, 10).Value = 835821
$ws.Cells.Item(26, 11).Value = 835820
$ws.Cells.Item(26, 12).Value = 835821
$ws.Cells.Item(26, 13).Value = 1156190

$ws.Cells.Item(27, 4).Value = 139
$ws.Cells.Item(27, 5).Value = 131
$ws.Cells.Item(27, 6).Value = 151
$ws.Cells.Item(27, 7).Value = 137
$ws.Cells.Item(27, 8).Value = 222
$ws.Cells.Item(27, 9).Value = 235
$ws.Cells.Item(27, 10).Value = 256
$ws.Cells.Item(27, 11).Value = 712
$ws.Cells.Item(27, 12).Value = 351
$ws.Cells.Item(27, 13).Value = 487
